# Add a single bottom/top border line to the 4 cells of the small
# one-column "(data)" / "(sudarymo vieta)" signature table, so the
# placeholder text has an underline above it.
#
# WdBorderType: wdBorderTop = -1, wdBorderBottom = -3
# WdLineStyle : wdLineStyleSingle = 1
# WdColor     : wdColorAutomatic = -16777216

$d = $word.ActiveDocument
$t = $d.Tables.Item(3)

function Set-CellBorder($cell, $borderId) {
    $border = $cell.Borders.Item($borderId)
    $border.LineStyle = 1
    $border.LineWidth = 2
    $border.Color = -16777216
    if ($borderId -eq -1) {
        $cell.Borders.DistanceFromTop = 0
    } else {
        $cell.Borders.DistanceFromBottom = 0
    }
}

# Row 1 (blank cell above "(data)") -> bottom border
Set-CellBorder $t.Cell(1, 1) -3

# Row 2 "(data)" -> top border
Set-CellBorder $t.Cell(2, 1) -1

# Row 3 (blank cell above "(sudarymo vieta)") -> bottom border
Set-CellBorder $t.Cell(3, 1) -3

# Row 4 "(sudarymo vieta)" -> top border
Set-CellBorder $t.Cell(4, 1) -1

Write-Output "borders applied"
